# Adds 20 new batal/meaning/source rows (rows 72-91) to Sheet1, pulling
# their text from the shared-strings additions in the target diff, and
# extends the sheet a further four blank styled rows (92-95) to match the
# new used range, mirroring an Excel "type a new row, tab, type, enter" pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the centered style ("s=2" in the sheet XML) used by every other
# data row, across the whole newly-touched range (including the trailing
# still-empty rows) in one shot.
$ws.Range("A72:C95").HorizontalAlignment = -4108

$sourceText = 'لال ءُ یاقوت'

$ws.Range("B72").Value = 'مرد پہ کار ءُ جہد، نیکہ دستاں بہ بند بہ بوزان پہ زان'
$ws.Range("A72").Value = 'بیگار مہ بو، بیکار مہ بو۔'
$ws.Range("C72").Value = $sourceText

$ws.Range("B73").Value = 'بیمار ءِ حال پُرسی خداءِ حال پُرسی اِنت، سواب اِنت، بیمار ءِ کرّا بیمار ءَ دلبڈی دیاں، آئی دل ءَ کمزور نہ کن انت'
$ws.Range("A73").Value = 'بیمار ءِ دل ءَ دور مہ دئے۔'
$ws.Range("C73").Value = $sourceText

$ws.Range("B74").Value = 'کس نہ زانت کہ کجام وہد ءَ بیماری کاینت، آ چہ اللہ ءِ نیمگ ءَ '
$ws.Range("A74").Value = 'بیماری کسی انتظار ءَ نہ کنت۔'
$ws.Range("C74").Value = $sourceText

$ws.Range("B75").Value = 'آ چیز کہ کسان اِنت چہ کسانی ءَ آئی سرا بِہ گر، بیماری بیت یا کہ زہگے بیت، دست ءَ کہ روت گڑا پشومانی کار نہ دنت'
$ws.Range("A75").Value = 'بیماری ءُ غم ءَ دیر مہ دار۔'
$ws.Range("C75").Value = $sourceText

$ws.Range("A76").Value = 'بیماری ءُ کار کسی لحاظ ءَ نہ کنت۔'
$ws.Range("B76").Value = 'بیماری چہ خدائی نیمگ ءَ کاینت، کار چہ خدا بنت، اے کسی لحاظ ءَ نہ کن انت۔'
$ws.Range("C76").Value = $sourceText

$ws.Range("B77").Value = 'ہمک کار ءِ سرا ہما وہد امءَ بہ گر'
$ws.Range("A77").Value = 'بیماری ءَ بے درمان مہ کن۔'
$ws.Range("C77").Value = $sourceText

$ws.Range("B78").Value = 'کسے ءَ را کہ براس نیست بے وس، بے کس نا وس بیت، چنکس کہ مالدار بہ بئے، بلے بے براسی‌گران اِنت'
$ws.Range("A78").Value = 'بے براسی نا وسی، بے کسی، بے وسی۔'
$ws.Range("C78").Value = $sourceText

$ws.Range("B79").Value = 'زہر ءُ ترندیں گپ دیوان ءَ حراب‌کنت'
$ws.Range("A79").Value = 'بےبراہ ایں گپ مجلس ءَ مرداکنت۔'
$ws.Range("C79").Value = $sourceText

$ws.Range("B80").Value = 'بے ننگیں مردم ءَ راکسی غیرت نہ گیپت'
$ws.Range("A80").Value = 'بےپِس ءَ نہ پیرک داریت نہ پِس۔'
$ws.Range("C80").Value = $sourceText

$ws.Range("B81").Value = 'ناتپاکی وت بربادی ے، کسے کہ  ناتپاک بوت گڑا آباد نہ بیت'
$ws.Range("A81").Value = 'بے تپاکی ہلاکی۔'
$ws.Range("C81").Value = $sourceText

$ws.Range("B82").Value = 'جان دزّ ءُ جان بڈّ لنگڑ بیت'
$ws.Range("A82").Value = 'بے جان بے نان بیت۔'
$ws.Range("C82").Value = $sourceText

$ws.Range("B83").Value = 'آکہ بے دین اِنت آئی رزق ءِ تہابرکت مان نہ بیت'
$ws.Range("A83").Value = 'بے دین ءِ رزق برکت نہ کنت۔'
$ws.Range("C83").Value = $sourceText

$ws.Range("B84").Value = 'بے زرّ ءَ راگرّ اِنت ، کسّے وتی نہ کنت'
$ws.Range("A84").Value = 'بے زَرّی گَرّے۔'
$ws.Range("C84").Value = $sourceText

$ws.Range("B85").Value = 'آ کہ بے سما اِنت، آئی ءَ را پہ گُشگ ءَ ہچ نہ بیت'
$ws.Range("A85").Value = 'بے سُد پہ پنت ءُ نصحیت ءَ سُد نہ کنت۔'
$ws.Range("C85").Value = $sourceText

$ws.Range("B86").Value = 'آ شہر کہ گوں تئی تالہ ءَ نہ ٹئیت گڑا لَڈّ ءُ بار بہ کن'
$ws.Range("A86").Value = 'بے سریں شہر ءَ مہ نند۔'
$ws.Range("C86").Value = $sourceText

$ws.Range("B87").Value = 'آ کار ءُ چیز کہ تاوان دینت آ یانی سر امہلہ بہ گر کہ تاوان دنت، آزہگے بہ بیت یا کہ نادراہی ے'
$ws.Range("A87").Value = 'بے سریں کار ءِ سر ءَ بہ گر۔'
$ws.Range("C87").Value = $sourceText

$ws.Range("B88").Value = 'بے ننگ ءُ بے ضمیریں مردم ءَ نوکر مہ کن'
$ws.Range("A88").Value = 'بے سریں مرد ءَ نوکر مہ کن۔'
$ws.Range("C88").Value = $sourceText

$ws.Range("B89").Value = 'جاہ ئے کہ روئے توشگ بہ زور، بے سلاح ءُ شور ءَ، بے سرپدی ءُ زانتکاری ءَ حبر مہ کن'
$ws.Range("A89").Value = 'بے سِلاح ءَ سپر مہ کن، بے صلاح ءَ حبر مہ کن۔'
$ws.Range("C89").Value = $sourceText

$ws.Range("B90").Value = 'ہر چیز ءِ تہا انصاپ بہ بیت، بے تور ءَ تور مہ کن وتی قیامت ءَ کوار مہ کن'
$ws.Range("A90").Value = 'بے شاہیم ءَ تول نہ بیت۔'
$ws.Range("C90").Value = $sourceText

$ws.Range("B91").Value = 'کول ءُ کرار گوں شوق ءَ بیت، گوں پہکیں مہر ءَ بیت'
$ws.Range("A91").Value = 'بے شوق ءَ کول نہ بیت۔'
$ws.Range("C91").Value = $sourceText

# Leave the cursor/selection where the author ended up after the last entry.
$ws.Range("A91").Select()
